$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3576.4
$ws.Range("I112").Value = 821.3333
$ws.Range("J112").Value = 4757.143
$ws.Range("K112").Value = 2463.9999
$ws.Range("L112").Value = 14271.429
$ws.Range("M112").Value = -1355.9999
$ws.Range("N112").Value = -16487.429

$ws.Range("H113").Value = 2086.6667
$ws.Range("I113").Value = 1831.4286
$ws.Range("K113").Value = 1831.4286
$ws.Range("M113").Value = 1422.5714

$ws.Range("H116").Value = 2982.8262
$ws.Range("I116").Value = 2686.0715
$ws.Range("J116").Value = 3444.4443
$ws.Range("K116").Value = 2686.0715
$ws.Range("L116").Value = 3444.4443
$ws.Range("M116").Value = 755.9285
$ws.Range("N116").Value = -10328.4443

$ws.Range("H129").Value = 1029.8646
$ws.Range("I129").Value = 1416.8182
$ws.Range("J129").Value = 979.7882
$ws.Range("K129").Value = 4250.4546
$ws.Range("L129").Value = 2939.3646
$ws.Range("M129").Value = 749.5454
$ws.Range("N129").Value = -12939.3646

$ws.Range("H132").Value = 43904.695
$ws.Range("I132").Value = 6892.3335
$ws.Range("K132").Value = 20677.0005
$ws.Range("M132").Value = -18147.0005

$ws.Range("H135").Value = 14706615
$ws.Range("I135").Value = 762.5172
$ws.Range("J135").Value = 100000560
$ws.Range("K135").Value = 6862.6548
$ws.Range("L135").Value = 900005040
$ws.Range("M135").Value = -4327.6548
$ws.Range("N135").Value = -900010110

$ws.Range("H137").Value = 4231.788
$ws.Range("I137").Value = 903.8946999999999
$ws.Range("K137").Value = 2711.6841
$ws.Range("M137").Value = -161.6840999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10872.632
$ws.Range("I32").Value = 9677.423000000001
$ws.Range("K32").Value = 9677.423000000001
$ws.Range("M32").Value = -9390.423000000001

$ws.Range("H61").Value = 2500.5625
$ws.Range("I61").Value = 1883.6364
$ws.Range("J61").Value = 3857.8
$ws.Range("K61").Value = 1883.6364
$ws.Range("L61").Value = 3857.8
$ws.Range("M61").Value = -1671.6364
$ws.Range("N61").Value = -4281.8

$ws.Range("H122").Value = 1446.6666
$ws.Range("I122").Value = 1446.6666
$ws.Range("K122").Value = 4339.9998
$ws.Range("M122").Value = -1889.9998

$ws.Range("H132").Value = 29414060
$ws.Range("I132").Value = 41667920
$ws.Range("K132").Value = 125003760
$ws.Range("M132").Value = -125001230

$ws.Range("H136").Value = 2500.5625
$ws.Range("I136").Value = 1883.6364
$ws.Range("J136").Value = 3857.8
$ws.Range("K136").Value = 5650.9092
$ws.Range("L136").Value = 11573.4
$ws.Range("M136").Value = -3100.9092
$ws.Range("N136").Value = -16673.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 105
$ws.Range("I5").Value = 105
$ws.Range("K5").Value = 105
$ws.Range("M5").Value = 8

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H134").Value = 2846.2134
$ws.Range("I134").Value = 2262.3513
$ws.Range("J134").Value = 3414.7104
$ws.Range("K134").Value = 6787.053899999999
$ws.Range("L134").Value = 10244.1312
$ws.Range("M134").Value = -4252.053899999999
$ws.Range("N134").Value = -15314.1312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 873.7
$ws.Range("I5").Value = 244.2
$ws.Range("J5").Value = 1503.2
$ws.Range("K5").Value = 244.2
$ws.Range("L5").Value = 1503.2
$ws.Range("M5").Value = -132.2
$ws.Range("N5").Value = -1727.2

$ws.Range("H10").Value = 1000000
$ws.Range("I10").Value = 1000000
$ws.Range("K10").Value = 1000000
$ws.Range("M10").Value = -999861

$ws.Range("H12").Value = 16674.916
$ws.Range("I12").Value = 1700
$ws.Range("J12").Value = 21666.555
$ws.Range("K12").Value = 1700
$ws.Range("L12").Value = 21666.555
$ws.Range("M12").Value = -1530
$ws.Range("N12").Value = -22006.555

$ws.Range("H14").Value = 40000
$ws.Range("J14").Value = 40000
$ws.Range("L14").Value = 40000
$ws.Range("N14").Value = -40340

$ws.Range("H25").Value = 33328
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 39593.6
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 39593.6
$ws.Range("M25").Value = -1826
$ws.Range("N25").Value = -39941.6

$ws.Range("H31").Value = 4046.85
$ws.Range("I31").Value = 1743.3
$ws.Range("K31").Value = 1743.3
$ws.Range("M31").Value = -1448.3

$ws.Range("H34").Value = 4046.85
$ws.Range("I34").Value = 1743.3
$ws.Range("K34").Value = 1743.3
$ws.Range("M34").Value = -1541.3

$ws.Range("H58").Value = 1779.52
$ws.Range("I58").Value = 1370.1578
$ws.Range("J58").Value = 3075.8333
$ws.Range("K58").Value = 1370.1578
$ws.Range("L58").Value = 3075.8333
$ws.Range("M58").Value = -1167.1578
$ws.Range("N58").Value = -3481.8333

$ws.Range("H132").Value = 144620.7
$ws.Range("I132").Value = 2199
$ws.Range("J132").Value = 180226.12
$ws.Range("K132").Value = 6597
$ws.Range("L132").Value = 540678.36
$ws.Range("M132").Value = -4067
$ws.Range("N132").Value = -545738.36

$ws.Range("H136").Value = 1779.52
$ws.Range("I136").Value = 1370.1578
$ws.Range("J136").Value = 3075.8333
$ws.Range("K136").Value = 4110.4734
$ws.Range("L136").Value = 9227.499899999999
$ws.Range("M136").Value = -1560.4734
$ws.Range("N136").Value = -14327.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 13889337
$ws.Range("J2").Value = 31250946
$ws.Range("L2").Value = 187505676
$ws.Range("N2").Value = -187505902

$ws.Range("H80").Value = 83521920
$ws.Range("I80").Value = 667666.7
$ws.Range("J80").Value = 111140000
$ws.Range("K80").Value = 2003000.1
$ws.Range("L80").Value = 333420000
$ws.Range("M80").Value = -2002064.1
$ws.Range("N80").Value = -333421872

$ws.Range("H83").Value = 83521920
$ws.Range("I83").Value = 667666.7
$ws.Range("J83").Value = 111140000
$ws.Range("K83").Value = 6009000.3
$ws.Range("L83").Value = 1000260000
$ws.Range("M83").Value = -6004320.3
$ws.Range("N83").Value = -1000269360

$ws.Range("H117").Value = 1929
$ws.Range("I117").Value = 315
$ws.Range("J117").Value = 2390.1428
$ws.Range("K117").Value = 945
$ws.Range("L117").Value = 7170.428400000001
$ws.Range("M117").Value = 2497
$ws.Range("N117").Value = -14054.4284

$ws.Range("H124").Value = 1993.6857
$ws.Range("I124").Value = 500
$ws.Range("J124").Value = 2037.6177
$ws.Range("K124").Value = 1500
$ws.Range("L124").Value = 6112.8531
$ws.Range("M124").Value = 3410
$ws.Range("N124").Value = -15932.8531

$ws.Range("H129").Value = 159486.27
$ws.Range("I129").Value = 500881.16
$ws.Range("J129").Value = 1919.3846
$ws.Range("K129").Value = 1502643.48
$ws.Range("L129").Value = 5758.1538
$ws.Range("M129").Value = -1497643.48
$ws.Range("N129").Value = -15758.1538

$ws.Range("H131").Value = 877.76
$ws.Range("J131").Value = 881.57574
$ws.Range("L131").Value = 2644.72722
$ws.Range("N131").Value = -12724.72722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 100
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -324

$ws.Range("H5").Value = 15702.5
$ws.Range("J5").Value = 15702.5
$ws.Range("L5").Value = 15702.5
$ws.Range("N5").Value = -15926.5

$ws.Range("H6").Value = 28000
$ws.Range("J6").Value = 28000
$ws.Range("L6").Value = 28000
$ws.Range("N6").Value = -28226

$ws.Range("H9").Value = 50605.332
$ws.Range("I9").Value = 11800
$ws.Range("K9").Value = 11800
$ws.Range("M9").Value = -11630

$ws.Range("H16").Value = 28000
$ws.Range("J16").Value = 28000
$ws.Range("L16").Value = 28000
$ws.Range("N16").Value = -28500

$ws.Range("H134").Value = 28500
$ws.Range("J134").Value = 28500
$ws.Range("L134").Value = 85500
$ws.Range("N134").Value = -90570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3139
$ws.Range("I9").Value = 228
$ws.Range("J9").Value = 6777.75
$ws.Range("K9").Value = 228
$ws.Range("L9").Value = 6777.75
$ws.Range("M9").Value = -4
$ws.Range("N9").Value = -7225.75

$ws.Range("H136").Value = 2013.6072
$ws.Range("I136").Value = 1512.9546
$ws.Range("J136").Value = 3849.3333
$ws.Range("K136").Value = 4538.8638
$ws.Range("L136").Value = 11547.9999
$ws.Range("M136").Value = -1988.8638
$ws.Range("N136").Value = -16647.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 10000935
$ws.Range("I107").Value = 335.7143
$ws.Range("J107").Value = 33335666
$ws.Range("K107").Value = 1007.1429
$ws.Range("L107").Value = 100006998
$ws.Range("M107").Value = 912.8571000000001
$ws.Range("N107").Value = -100010838

$ws.Range("H126").Value = 2487.2666
$ws.Range("I126").Value = 1832.7142
$ws.Range("J126").Value = 3060
$ws.Range("K126").Value = 5498.142599999999
$ws.Range("L126").Value = 9180
$ws.Range("M126").Value = -3028.142599999999
$ws.Range("N126").Value = -14120

$ws.Range("H132").Value = 1595.3572
$ws.Range("I132").Value = 1215.129
$ws.Range("J132").Value = 2666.9092
$ws.Range("K132").Value = 3645.387
$ws.Range("L132").Value = 8000.7276
$ws.Range("M132").Value = -1115.387
$ws.Range("N132").Value = -13060.7276

$ws.Range("H136").Value = 436320.47
$ws.Range("I136").Value = 667838.25
$ws.Range("J136").Value = 2224.625
$ws.Range("K136").Value = 2003514.75
$ws.Range("L136").Value = 6673.875
$ws.Range("M136").Value = -2000964.75
$ws.Range("N136").Value = -11773.875
